$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.950.11"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "3.052.38"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.91%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.049.81"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.37%  "
$ws.Range("E10").Value = "  +6.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.478"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.19%  "
$ws.Range("E13").Value = "  +6.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "3.537.08"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").Value = "64.002.15"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("D17").Value = "3.056.97"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.19%  "
$ws.Range("E22").Value = "  +4.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +14.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +3.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.74%  "
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("E33").Value = "  +5.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("E35").Value = "  +7.19%  "
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0405"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "441.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.66%  "
$ws.Range("D41").Value = "2.969.32"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.36%  "
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("E44").Value = "  +4.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.82%  "
$ws.Range("E46").Value = "  +9.35%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  +4.89%  "
$ws.Range("E49").Value = "  +6.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.55%  "
